# Fill in the Planning table on the "Iteracion4" sheet with the missing
# estimation values for "Ana" (col B), "Santi Lopez" (col G) and the
# final "Estimación" (col J), then update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteracion4")
$ws.Activate()

# Row 4 - "Mejora: Cambio de logo"
$ws.Range("B4").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("J4").Value = 3

# Row 5 - "Mejora: Agregar modo nocturno"
$ws.Range("B5").Value = 8
$ws.Range("G5").Value = 8
$ws.Range("J5").Value = 13

# Row 6 - "Mejora: Agregar idiomas"
$ws.Range("B6").Value = 5
$ws.Range("G6").Value = 8
$ws.Range("J6").Value = 8

# Row 7 - "Agregar chat"
$ws.Range("B7").Value = 8
$ws.Range("G7").Value = 8
$ws.Range("J7").Value = 8

# Row 8 - "Mejoras para los favoritos"
$ws.Range("B8").Value = 3
$ws.Range("G8").Value = 3
$ws.Range("J8").Value = 3

# Update the selected cell shown in the saved view.
$ws.Range("F19").Select()
